$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new error code row (row 41)
$ws.Range("A41").Value = 412039
$ws.Range("B41").Value = "SerializationException - An error occurred during serialization. This is probably due to a JSON payload being malformed."

# Reflect the cursor position left after entering the new row's data
$ws.Range("B42").Select()
